$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 54
$ws_ALC.Range("H54").Value = 136085.62
$ws_ALC.Range("I54").Value = 204514.8
$ws_ALC.Range("J54").Value = 22037
$ws_ALC.Range("K54").Value = 204514.8
$ws_ALC.Range("L54").Value = 22037
$ws_ALC.Range("M54").Value = -204028.8
$ws_ALC.Range("N54").Value = -23009

# ALC row 64
$ws_ALC.Range("H64").Value = 8335.474
$ws_ALC.Range("J64").Value = 8971.214
$ws_ALC.Range("L64").Value = 8971.214
$ws_ALC.Range("N64").Value = -9467.214

# ALC row 67
$ws_ALC.Range("H67").Value = 8335.474
$ws_ALC.Range("J67").Value = 8971.214
$ws_ALC.Range("L67").Value = 8971.214
$ws_ALC.Range("N67").Value = -10687.214

# ALC row 70
$ws_ALC.Range("H70").Value = 10125.167
$ws_ALC.Range("I70").Value = 7416.3335
$ws_ALC.Range("K70").Value = 22249.0005
$ws_ALC.Range("M70").Value = -21979.0005

# ALC row 73
$ws_ALC.Range("H73").Value = 10125.167
$ws_ALC.Range("I73").Value = 7416.3335
$ws_ALC.Range("K73").Value = 22249.0005
$ws_ALC.Range("M73").Value = -21313.0005

# ALC row 98
$ws_ALC.Range("H98").Value = 214955.64
$ws_ALC.Range("I98").Value = 780.3333
$ws_ALC.Range("J98").Value = 857481.5600000001
$ws_ALC.Range("K98").Value = 780.3333
$ws_ALC.Range("L98").Value = 857481.5600000001
$ws_ALC.Range("M98").Value = 717.6667
$ws_ALC.Range("N98").Value = -860477.5600000001

# ALC row 107
$ws_ALC.Range("H107").Value = 356.6
$ws_ALC.Range("I107").Value = 245.75
$ws_ALC.Range("K107").Value = 245.75
$ws_ALC.Range("M107").Value = 1674.25

# ALC row 122
$ws_ALC.Range("H122").Value = 214955.64
$ws_ALC.Range("I122").Value = 780.3333
$ws_ALC.Range("J122").Value = 857481.5600000001
$ws_ALC.Range("K122").Value = 2340.9999
$ws_ALC.Range("L122").Value = 2572444.68
$ws_ALC.Range("M122").Value = 109.0001000000002
$ws_ALC.Range("N122").Value = -2577344.68

# ALC row 141
$ws_ALC.Range("H141").Value = 3848
$ws_ALC.Range("I141").Value = 2746.6667
$ws_ALC.Range("K141").Value = 8240.000100000001
$ws_ALC.Range("M141").Value = -3060.000100000001

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws_ARM.Range("H32").Value = 3396.7754
$ws_ARM.Range("I32").Value = 2748.739
$ws_ARM.Range("K32").Value = 2748.739
$ws_ARM.Range("M32").Value = -2461.739

# ARM row 61
$ws_ARM.Range("H61").Value = 4546
$ws_ARM.Range("I61").Value = 3311.3635
$ws_ARM.Range("K61").Value = 3311.3635
$ws_ARM.Range("M61").Value = -3099.3635

# ARM row 136
$ws_ARM.Range("H136").Value = 4546
$ws_ARM.Range("I136").Value = 3311.3635
$ws_ARM.Range("K136").Value = 9934.0905
$ws_ARM.Range("M136").Value = -7384.0905

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws_BSM.Range("H105").Value = 15756.35
$ws_BSM.Range("J105").Value = 12771.429
$ws_BSM.Range("L105").Value = 12771.429
$ws_BSM.Range("N105").Value = -16265.429

# BSM row 107
$ws_BSM.Range("H107").Value = 2373.375
$ws_BSM.Range("J107").Value = 2569.5715
$ws_BSM.Range("L107").Value = 2569.5715
$ws_BSM.Range("N107").Value = -6409.5715

# BSM row 132
$ws_BSM.Range("H132").Value = 54995.855
$ws_BSM.Range("J132").Value = 54995.855
$ws_BSM.Range("L132").Value = 54995.855
$ws_BSM.Range("N132").Value = -65115.855

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws_CRP.Range("H22").Value = 2039.8889
$ws_CRP.Range("I22").Value = 476.66666
$ws_CRP.Range("K22").Value = 476.66666
$ws_CRP.Range("M22").Value = -126.66666

# CRP row 74
$ws_CRP.Range("H74").Value = 0
$ws_CRP.Range("J74").Value = 0
$ws_CRP.Range("N74").Value = 0
$ws_CRP.Range("L74").Value = ""

# CRP row 77
$ws_CRP.Range("H77").Value = 0
$ws_CRP.Range("J77").Value = 0
$ws_CRP.Range("N77").Value = 0
$ws_CRP.Range("L77").Value = ""

# CRP row 134
$ws_CRP.Range("H134").Value = 3733.5557
$ws_CRP.Range("I134").Value = 2636.875
$ws_CRP.Range("K134").Value = 7910.625
$ws_CRP.Range("M134").Value = -5375.625

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 2
$ws_CUL.Range("H2").Value = 174.17392
$ws_CUL.Range("I2").Value = 126.07143
$ws_CUL.Range("J2").Value = 249
$ws_CUL.Range("K2").Value = 756.42858
$ws_CUL.Range("L2").Value = 1494
$ws_CUL.Range("M2").Value = -643.42858
$ws_CUL.Range("N2").Value = -1720

# CUL row 16
$ws_CUL.Range("H16").Value = 2855.7144
$ws_CUL.Range("I16").Value = 2495
$ws_CUL.Range("K16").Value = 7485
$ws_CUL.Range("M16").Value = -7312

# CUL row 37
$ws_CUL.Range("H37").Value = 197222.56
$ws_CUL.Range("J37").Value = 197222.56
$ws_CUL.Range("L37").Value = 591667.6799999999
$ws_CUL.Range("N37").Value = -591891.6799999999

# CUL row 38
$ws_CUL.Range("H38").Value = 77.5
$ws_CUL.Range("I38").Value = 73
$ws_CUL.Range("K38").Value = 219
$ws_CUL.Range("M38").Value = 128

# CUL row 60
$ws_CUL.Range("H60").Value = 41667940
$ws_CUL.Range("I60").Value = 55555920
$ws_CUL.Range("J60").Value = 4000
$ws_CUL.Range("K60").Value = 166667760
$ws_CUL.Range("L60").Value = 12000
$ws_CUL.Range("M60").Value = -166667509
$ws_CUL.Range("N60").Value = -12502

# CUL row 63
$ws_CUL.Range("H63").Value = 10000
$ws_CUL.Range("I63").Value = 0
$ws_CUL.Range("J63").Value = 10000
$ws_CUL.Range("K63").Value = 0
$ws_CUL.Range("M63").Value = 30000
$ws_CUL.Range("N63").Value = -31498
$ws_CUL.Range("L63").Value = ""

# CUL row 66
$ws_CUL.Range("H66").Value = 10000
$ws_CUL.Range("I66").Value = 0
$ws_CUL.Range("J66").Value = 10000
$ws_CUL.Range("K66").Value = 0
$ws_CUL.Range("M66").Value = 90000
$ws_CUL.Range("N66").Value = -97488
$ws_CUL.Range("L66").Value = ""

# CUL row 75
$ws_CUL.Range("H75").Value = 111116280
$ws_CUL.Range("I75").Value = 333335680
$ws_CUL.Range("J75").Value = 6586.3335
$ws_CUL.Range("K75").Value = 1000007040
$ws_CUL.Range("L75").Value = 19759.0005
$ws_CUL.Range("M75").Value = -1000006042
$ws_CUL.Range("N75").Value = -21755.0005

# CUL row 76
$ws_CUL.Range("H76").Value = 0
$ws_CUL.Range("I76").Value = 0
$ws_CUL.Range("K76").Value = 0
$ws_CUL.Range("M76").Value = ""

# CUL row 78
$ws_CUL.Range("H78").Value = 111116280
$ws_CUL.Range("I78").Value = 333335680
$ws_CUL.Range("J78").Value = 6586.3335
$ws_CUL.Range("K78").Value = 3000021120
$ws_CUL.Range("L78").Value = 59277.0015
$ws_CUL.Range("M78").Value = -3000016128
$ws_CUL.Range("N78").Value = -69261.0015

# CUL row 79
$ws_CUL.Range("H79").Value = 0
$ws_CUL.Range("I79").Value = 0
$ws_CUL.Range("K79").Value = 0
$ws_CUL.Range("M79").Value = ""

# CUL row 81
$ws_CUL.Range("H81").Value = 5434.6665
$ws_CUL.Range("I81").Value = 1980.4
$ws_CUL.Range("K81").Value = 5941.200000000001
$ws_CUL.Range("M81").Value = -4818.200000000001

# CUL row 84
$ws_CUL.Range("H84").Value = 5434.6665
$ws_CUL.Range("I84").Value = 1980.4
$ws_CUL.Range("K84").Value = 17823.6
$ws_CUL.Range("M84").Value = -12207.6

# CUL row 86
$ws_CUL.Range("H86").Value = 3057.6924
$ws_CUL.Range("I86").Value = 2918.5
$ws_CUL.Range("K86").Value = 8755.5
$ws_CUL.Range("M86").Value = -7569.5

# CUL row 88
$ws_CUL.Range("I88").Value = 20000
$ws_CUL.Range("J88").Value = 0
$ws_CUL.Range("K88").Value = 60000
$ws_CUL.Range("N88").Value = 0
$ws_CUL.Range("L88").Value = ""
$ws_CUL.Range("M88").Value = -59572

# CUL row 89
$ws_CUL.Range("H89").Value = 3057.6924
$ws_CUL.Range("I89").Value = 2918.5
$ws_CUL.Range("K89").Value = 26266.5
$ws_CUL.Range("M89").Value = -20338.5

# CUL row 91
$ws_CUL.Range("I91").Value = 20000
$ws_CUL.Range("J91").Value = 0
$ws_CUL.Range("K91").Value = 60000
$ws_CUL.Range("N91").Value = 0
$ws_CUL.Range("L91").Value = ""
$ws_CUL.Range("M91").Value = -58518

# CUL row 136
$ws_CUL.Range("H136").Value = 30306384
$ws_CUL.Range("I136").Value = 37040024
$ws_CUL.Range("K136").Value = 111120072
$ws_CUL.Range("M136").Value = -111114972

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws_GSM.Range("H102").Value = 1712.228
$ws_GSM.Range("I102").Value = 1259.5
$ws_GSM.Range("K102").Value = 1259.5
$ws_GSM.Range("M102").Value = 362.5

# GSM row 122
$ws_GSM.Range("H122").Value = 17832.2
$ws_GSM.Range("I122").Value = 21267.715
$ws_GSM.Range("J122").Value = 9816
$ws_GSM.Range("K122").Value = 63803.145
$ws_GSM.Range("L122").Value = 29448
$ws_GSM.Range("M122").Value = -61353.145
$ws_GSM.Range("N122").Value = -34348

# GSM row 132
$ws_GSM.Range("H132").Value = 9321.416999999999
$ws_GSM.Range("I132").Value = 8440.272000000001
$ws_GSM.Range("K132").Value = 25320.816
$ws_GSM.Range("M132").Value = -22790.816

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws_LTW.Range("H22").Value = 3895.75
$ws_LTW.Range("I22").Value = 1204.8823
$ws_LTW.Range("K22").Value = 1204.8823
$ws_LTW.Range("M22").Value = -909.8823

# LTW row 27
$ws_LTW.Range("H27").Value = 3895.75
$ws_LTW.Range("I27").Value = 1204.8823
$ws_LTW.Range("K27").Value = 1204.8823
$ws_LTW.Range("M27").Value = -1097.8823

# LTW row 40
$ws_LTW.Range("H40").Value = 15975.143
$ws_LTW.Range("I40").Value = 18327.777
$ws_LTW.Range("J40").Value = 11740.4
$ws_LTW.Range("K40").Value = 18327.777
$ws_LTW.Range("L40").Value = 11740.4
$ws_LTW.Range("M40").Value = -18191.777
$ws_LTW.Range("N40").Value = -12012.4

# LTW row 61
$ws_LTW.Range("H61").Value = 15034.617
$ws_LTW.Range("I61").Value = 17313.74
$ws_LTW.Range("J61").Value = 6243.7144
$ws_LTW.Range("K61").Value = 17313.74
$ws_LTW.Range("L61").Value = 6243.7144
$ws_LTW.Range("M61").Value = -17111.74
$ws_LTW.Range("N61").Value = -6647.7144

# LTW row 113
$ws_LTW.Range("H113").Value = 15034.617
$ws_LTW.Range("I113").Value = 17313.74
$ws_LTW.Range("J113").Value = 6243.7144
$ws_LTW.Range("K113").Value = 17313.74
$ws_LTW.Range("L113").Value = 6243.7144
$ws_LTW.Range("M113").Value = -15143.74
$ws_LTW.Range("N113").Value = -10583.7144

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 30
$ws_WVR.Range("H30").Value = 10955
$ws_WVR.Range("I30").Value = 10955
$ws_WVR.Range("J30").Value = 0
$ws_WVR.Range("K30").Value = 10955
$ws_WVR.Range("N30").Value = 0
$ws_WVR.Range("L30").Value = ""
$ws_WVR.Range("M30").Value = -10848

# WVR row 46
$ws_WVR.Range("H46").Value = 0
$ws_WVR.Range("J46").Value = 0
$ws_WVR.Range("N46").Value = 0
$ws_WVR.Range("L46").Value = ""

# WVR row 81
$ws_WVR.Range("H81").Value = 3569.8
$ws_WVR.Range("I81").Value = 1334.5
$ws_WVR.Range("J81").Value = 8040.4
$ws_WVR.Range("K81").Value = 2669
$ws_WVR.Range("L81").Value = 16080.8
$ws_WVR.Range("M81").Value = -1608
$ws_WVR.Range("N81").Value = -18202.8

# WVR row 84
$ws_WVR.Range("H84").Value = 3569.8
$ws_WVR.Range("I84").Value = 1334.5
$ws_WVR.Range("J84").Value = 8040.4
$ws_WVR.Range("K84").Value = 13345
$ws_WVR.Range("L84").Value = 80404
$ws_WVR.Range("M84").Value = -8041
$ws_WVR.Range("N84").Value = -91012

# WVR row 134
$ws_WVR.Range("H134").Value = 0
$ws_WVR.Range("J134").Value = 0
$ws_WVR.Range("N134").Value = 0
$ws_WVR.Range("L134").Value = ""

